$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (625.x264_s) now has a re-run result instead of "Run crashed" ---
# Un-merge the D10:E10 cell (it previously held a single merged "Run crashed" label)
$ws.Range("D10:E10").UnMerge()

# Copy the normal numeric-row formatting (borders/number format/font) from row 6
# onto row 10's D:F cells, since row 10 loses its special "crashed" styling.
$ws.Range("D6:F6").Copy()
$ws.Range("D10:F10").PasteSpecial(-4122)

# Fill in the new measured values for the re-run of 625.x264_s
$ws.Range("D10").Value = 4333
$ws.Range("E10").Value = 0.40699999999999997
$ws.Range("F10").Formula = "=D10/B10"

# --- Annotate the Base Score header with a footnote marker ---
$ws.Range("D15").Value = "Base Score**"

# --- Add the footnote itself below the existing notes ---
$ws.Range("A21").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = "**:"
$ws.Range("B25").Value = "Not accurate, as the 625.x264_s data is from a separate run"

# --- View state tweaks captured in the saved workbook ---
$excel.ActiveWindow.Zoom = 142
$ws.Range("C24").Select()
